$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("experiment_description")
$ws2 = $wb.Worksheets.Item("experiment_specification")
$ws3 = $wb.Worksheets.Item("run_description")
$ws4 = $wb.Worksheets.Item("run_specification")

# --- experiment_description (rows 20-24) & run_description (rows 17-23) ---
# Interleaved to reproduce authentic shared-string insertion order
$ws1.Range("A20").Value = 16
$ws1.Range("B20").Value = 'nytt'
$ws1.Range("C20").Value = 'base'
$ws1.Range("D20").Value = 1
$ws1.Range("E20").Value = 'none'
$ws1.Range("F20").Value = 'none'
$ws3.Range("A17").Value = 16
$ws3.Range("B17").Value = 'nytt prufa'
$ws1.Range("A21").Value = 2001
$ws1.Range("B21").Value = 'base model with date_first_symptoms instead of date_diagnosis'
$ws1.Range("C21").Value = 'symptoms'
$ws1.Range("D21").Value = '1;2;3'
$ws1.Range("E21").Value = 'none'
$ws1.Range("F21").Value = 'none'
$ws1.Range("G21").Value = 'Date of first symptoms'
$ws1.Range("H21").Value = 'Dagseting fyrstu einkenna'
$ws3.Range("A18").Value = 2001
$ws3.Range("B18").Value = 'date_symptoms'
$ws1.Range("A22").Value = 2004
$ws1.Range("B22").Value = 'LOS driven: 3 + age-point of diag splitting of transitions in Inpatient Ward, first_symptoms'
$ws1.Range("C22").Value = 'symptoms'
$ws1.Range("D22").Value = '1;4'
$ws1.Range("E22").Value = 'none'
$ws1.Range("F22").Value = 'none'
$ws1.Range("G22").Value = 'symptoms improved'
$ws1.Range("H22").Value = 'bætt líkan með symptoms'
$ws3.Range("A19").Value = 2004
$ws3.Range("B19").Value = '4 með date_symptoms'
$ws1.Range("A23").Value = 2006
$ws1.Range("B23").Value = 'Transition driven: 1 + hospital_less_than_14_days heuristic,first_symptoms'
$ws1.Range("C23").Value = 'symptoms'
$ws1.Range("D23").Value = '1;2;3;5'
$ws1.Range("E23").Value = 'none'
$ws1.Range("F23").Value = 'none'
$ws1.Range("G23").Value = 6
$ws1.Range("H23").Value = 6
$ws1.Range("A24").Value = 2007
$ws1.Range("B24").Value = 'LOS driven:3 + treatment constraints splitting of transitions  in Inpatient Ward,first_symptoms'
$ws1.Range("C24").Value = 'symptoms'
$ws1.Range("D24").Value = '1;4'
$ws1.Range("E24").Value = 'none'
$ws1.Range("F24").Value = 'none'
$ws1.Range("G24").Value = 7
$ws1.Range("H24").Value = 7
$ws3.Range("A20").Value = 6
$ws3.Range("B20").Value = '6 with date_symptoms'
$ws3.Range("A21").Value = 7
$ws3.Range("B21").Value = '7 with date_symptoms'
$ws3.Range("A22").Value = 17
$ws3.Range("B22").Value = 'Only improved model'
$ws3.Range("A23").Value = 18
$ws3.Range("B23").Value = 'Only ferguson wuhan'

# --- experiment_specification: add rows 65-79 (reuses existing strings only) ---
$ws2.Range("A65").Value = 16
$ws2.Range("B65").Value = 'home'
$ws2.Range("C65").Value = 'none'
$ws2.Range("D65").Value = 'age_simple'
$ws2.Range("E65").Value = 'none'
$ws2.Range("A66").Value = 16
$ws2.Range("B66").Value = 'inpatient_ward'
$ws2.Range("C66").Value = 'none'
$ws2.Range("D66").Value = 'none'
$ws2.Range("E66").Value = 'none'
$ws2.Range("A67").Value = 16
$ws2.Range("B67").Value = 'intensive_care_unit'
$ws2.Range("C67").Value = 'none'
$ws2.Range("D67").Value = 'none'
$ws2.Range("E67").Value = 'none'
$ws2.Range("A68").Value = 2001
$ws2.Range("B68").Value = 'home'
$ws2.Range("C68").Value = 'none'
$ws2.Range("D68").Value = 'age_simple'
$ws2.Range("E68").Value = 'age_simple'
$ws2.Range("A69").Value = 2001
$ws2.Range("B69").Value = 'inpatient_ward'
$ws2.Range("C69").Value = 'none'
$ws2.Range("D69").Value = 'age_simple'
$ws2.Range("E69").Value = 'none'
$ws2.Range("A70").Value = 2001
$ws2.Range("B70").Value = 'intensive_care_unit'
$ws2.Range("C70").Value = 'none'
$ws2.Range("D70").Value = 'age_simple'
$ws2.Range("E70").Value = 'none'
$ws2.Range("A71").Value = 2004
$ws2.Range("B71").Value = 'home'
$ws2.Range("C71").Value = 'length_of_stay_simple_two_weeks'
$ws2.Range("D71").Value = 'age_simple'
$ws2.Range("E71").Value = 'age_simple'
$ws2.Range("A72").Value = 2004
$ws2.Range("B72").Value = 'inpatient_ward'
$ws2.Range("C72").Value = 'none'
$ws2.Range("D72").Value = 'age_simple_point_of_diagnosis'
$ws2.Range("E72").Value = 'none'
$ws2.Range("A73").Value = 2004
$ws2.Range("B73").Value = 'intensive_care_unit'
$ws2.Range("C73").Value = 'none'
$ws2.Range("D73").Value = 'age_simple'
$ws2.Range("E73").Value = 'none'
$ws2.Range("A74").Value = 2006
$ws2.Range("B74").Value = 'home'
$ws2.Range("C74").Value = 'none'
$ws2.Range("D74").Value = 'age_simple'
$ws2.Range("E74").Value = 'age_simple'
$ws2.Range("A75").Value = 2006
$ws2.Range("B75").Value = 'inpatient_ward'
$ws2.Range("C75").Value = 'none'
$ws2.Range("D75").Value = 'age_simple'
$ws2.Range("E75").Value = 'none'
$ws2.Range("A76").Value = 2006
$ws2.Range("B76").Value = 'intensive_care_unit'
$ws2.Range("C76").Value = 'none'
$ws2.Range("D76").Value = 'age_simple'
$ws2.Range("E76").Value = 'none'
$ws2.Range("A77").Value = 2007
$ws2.Range("B77").Value = 'home'
$ws2.Range("C77").Value = 'length_of_stay_simple_two_weeks'
$ws2.Range("D77").Value = 'age_simple'
$ws2.Range("E77").Value = 'age_simple'
$ws2.Range("A78").Value = 2007
$ws2.Range("B78").Value = 'inpatient_ward'
$ws2.Range("C78").Value = 'none'
$ws2.Range("D78").Value = 'age_simple_intensive_care_unit_restriction'
$ws2.Range("E78").Value = 'none'
$ws2.Range("A79").Value = 2007
$ws2.Range("B79").Value = 'intensive_care_unit'
$ws2.Range("C79").Value = 'none'
$ws2.Range("D79").Value = 'age_simple'
$ws2.Range("E79").Value = 'none'

# --- run_specification: add rows 41-48 ---
$ws4.Range("A41").Value = 16
$ws4.Range("B41").Value = 16
$ws4.Range("A42").Value = 2001
$ws4.Range("B42").Value = 2001
$ws4.Range("A43").Value = 2004
$ws4.Range("B43").Value = 2004
$ws4.Range("A44").Value = 2006
$ws4.Range("B44").Value = 2006
$ws4.Range("A45").Value = 2006
$ws4.Range("B45").Value = 2007
$ws4.Range("A46").Value = 2007
$ws4.Range("B46").Value = 2007
$ws4.Range("A47").Value = 17
$ws4.Range("B47").Value = 4
$ws4.Range("A48").Value = 18
$ws4.Range("B48").Value = 10
